$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Value = "Leatitia"
$ws.Range("C5").Value = "Urban corporation"
$ws.Range("F5").Value = "urbain.calt@gmail.com"
$ws.Columns("B:C").AutoFit()
$ws.Columns("F:F").AutoFit()
Write-Output "done"
